$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Sat Aug 17 22:26:58 UTC 2024 with GitHub Actions

$ws.Range("D2").Value = "59.232.46"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "2.604.49"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'540.59"
$ws.Range("E5").Value = "  +3.85%  "
$ws.Range("D6").Value = "'141.58"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").Value = "'0.335"
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "3.061.40"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "59.169.85"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "'20.59"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "2.623.08"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "'341.90"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "'4.37"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").Value = "'10.14"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "'6.38"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'67.58"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("D28").Value = "0.0₃0758"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'1.70"
$ws.Range("E30").Value = "  +8.20%  "
$ws.Range("D31").Value = "'5.84"
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").Value = "'18.74"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").Value = "'149.44"
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").Value = "'37.16"
$ws.Range("E36").Value = "  +2.48%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "'0.836"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").Value = "'0.825"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "'274.37"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").Value = "'0.597"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "'0.0957"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "'0.0525"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.954.36"
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'18.59"
$ws.Range("E48").Value = "  +3.49%  "
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "'111.09"
$ws.Range("E51").Value = "  -1.14%  "

# Clear the auto-applied "stored as text" style so cells keep the workbook default formatting
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
